# paises.xlsx update: "Update countries & provincias Spain"
# - Refreshes the COVID-19 per-country figures for the countries whose case
#   counts changed, and re-applies the "sorted by Casos totales desc" order,
#   which shuffles a few country labels (Marruecos, Etiopia, Libano each
#   moved up one rank and displaced the country that used to sit above them).
# - Updates the "last refreshed" timestamp banner in A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Datos actualizados ..." timestamp banner ---
$ws.Range("A1").Value = "Datos actualizados a 2 de Agosto de 2020 a las 19:46"

# Row 4: Estados Unidos - updated figures
$ws.Range("B4").Value = 4789949
$ws.Range("C4").Value = 25631
$ws.Range("D4").Value = 2371733
$ws.Range("E4").Value = 2260044
$ws.Range("G4").Value = 274
$ws.Range("H4").Value = 158172

# Row 5: Brasil - updated figures
$ws.Range("B5").Value = 2711132
$ws.Range("C5").Value = 2256
$ws.Range("E5").Value = 733422
$ws.Range("G5").Value = 43
$ws.Range("H5").Value = 93659

# Row 6: India - updated figures
$ws.Range("B6").Value = 1804258
$ws.Range("C6").Value = 52339
$ws.Range("D6").Value = 1186981
$ws.Range("E6").Value = 579119
$ws.Range("G6").Value = 755
$ws.Range("H6").Value = 38158

# Row 11: Chile - updated figures
$ws.Range("B11").Value = 359731
$ws.Range("C11").Value = 2073
$ws.Range("D11").Value = 332411
$ws.Range("E11").Value = 17712
$ws.Range("G11").Value = 75
$ws.Range("H11").Value = 9608

# Row 20: Turquia - updated figures
$ws.Range("B20").Value = 232856
$ws.Range("C20").Value = 987
$ws.Range("D20").Value = 216494
$ws.Range("E20").Value = 10634
$ws.Range("G20").Value = 18
$ws.Range("H20").Value = 5728

# Row 36: Israel - updated figures
$ws.Range("B36").Value = 72584
$ws.Range("C36").Value = 366
$ws.Range("D36").Value = 45662
$ws.Range("E36").Value = 26386
$ws.Range("G36").Value = 10
$ws.Range("H36").Value = 536

# Row 62: Irlanda - updated figures
$ws.Range("B62").Value = 26162
$ws.Range("C62").Value = 53
$ws.Range("E62").Value = 1035

# Row 63: now "Marruecos" (was "Moldavia", displaced by re-sort)
$ws.Range("A63").Value = "Marruecos"
$ws.Range("B63").Value = 25537
$ws.Range("C63").Value = 522
$ws.Range("D63").Value = 18435
$ws.Range("E63").Value = 6720
$ws.Range("G63").Value = 15
$ws.Range("H63").Value = 382

# Row 64: now "Moldavia" (was "Uzbekistan", displaced by re-sort)
$ws.Range("A64").Value = "Moldavia"
$ws.Range("B64").Value = 25362
$ws.Range("C64").Value = 249
$ws.Range("D64").Value = 17816
$ws.Range("E64").Value = 6755
$ws.Range("G64").Value = 3
$ws.Range("H64").Value = 791

# Row 65: now "Uzbekistan" (was "Marruecos", displaced by re-sort)
$ws.Range("A65").Value = "Uzbekistan"
$ws.Range("B65").Value = 25336
$ws.Range("C65").Value = 553
$ws.Range("D65").Value = 15833
$ws.Range("E65").Value = 9352
$ws.Range("G65").Value = 4
$ws.Range("H65").Value = 151

# Row 70: now "Etiopia" (was "Costa Rica", displaced by re-sort)
$ws.Range("A70").Value = "Etiopia"
$ws.Range("B70").Value = 18706
$ws.Range("C70").Value = 707
$ws.Range("D70").Value = 7601
$ws.Range("E70").Value = 10795
$ws.Range("G70").Value = 26
$ws.Range("H70").Value = 310

# Row 71: now "Costa Rica" (was "Etiopia", displaced by re-sort)
$ws.Range("A71").Value = "Costa Rica"
$ws.Range("B71").Value = 18187
$ws.Range("D71").Value = 4531
$ws.Range("E71").Value = 13502
$ws.Range("H71").Value = 154

# Row 76: Costa de Marfil - updated figures
$ws.Range("B76").Value = 16182
$ws.Range("C76").Value = 73
$ws.Range("D76").Value = 11801
$ws.Range("E76").Value = 4279

# Row 102: now "Libano" (was "Guinea Ecuatorial", displaced by re-sort)
$ws.Range("A102").Value = "Libano"
$ws.Range("B102").Value = 4885
$ws.Range("C102").Value = 155
$ws.Range("D102").Value = 1795
$ws.Range("E102").Value = 3028
$ws.Range("G102").Value = 1
$ws.Range("H102").Value = 62

# Row 103: now "Guinea Ecuatorial" (was "Libano", displaced by re-sort)
$ws.Range("A103").Value = "Guinea Ecuatorial"
$ws.Range("B103").Value = 4821
$ws.Range("D103").Value = 2182
$ws.Range("E103").Value = 2556
$ws.Range("H103").Value = 83

# Row 107: Malaui - updated figures
$ws.Range("B107").Value = 4231
$ws.Range("C107").Value = 45
$ws.Range("D107").Value = 1919
$ws.Range("E107").Value = 2189
$ws.Range("G107").Value = 3
$ws.Range("H107").Value = 123

# Row 108: Maldivas - updated figures
$ws.Range("B108").Value = 4164
$ws.Range("C108").Value = 215
$ws.Range("D108").Value = 2643
$ws.Range("E108").Value = 1503
$ws.Range("G108").Value = 1
$ws.Range("H108").Value = 18

# Row 133: Sierra Leona - updated figures
$ws.Range("B133").Value = 1843
$ws.Range("C133").Value = 20
$ws.Range("D133").Value = 1375
$ws.Range("E133").Value = 401

# Row 136: Yemen - updated figures
$ws.Range("B136").Value = 1734
$ws.Range("C136").Value = 4
$ws.Range("E136").Value = 375
$ws.Range("G136").Value = 3
$ws.Range("H136").Value = 497

# Row 138: Tunez - updated figures
$ws.Range("B138").Value = 1561
$ws.Range("C138").Value = 9
$ws.Range("D138").Value = 1221
$ws.Range("E138").Value = 289

# Row 141: Jordania - updated figures
$ws.Range("B141").Value = 1213
$ws.Range("C141").Value = 5
$ws.Range("D141").Value = 1099

# Row 148: Niger - updated figures
$ws.Range("B148").Value = 1138
$ws.Range("C148").Value = 2
$ws.Range("E148").Value = 41

# Row 177: Islas Feroe - updated figures
$ws.Range("D177").Value = 192
$ws.Range("E177").Value = 33

# Row 185: Aruba - updated figures
$ws.Range("B185").Value = 122
$ws.Range("C185").Value = 1
$ws.Range("D185").Value = 111
$ws.Range("E185").Value = 8
